$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Healthcare Senior Golang Developer", "https://www.dice.com/job-detail/ae1b61df-4724-423c-9b66-1a37695d2877", "Remote or California", "Third Party", "Depends on Experience", "Pioneer Corporate Services Inc"),
    @("Java Backend Engineer", "https://www.dice.com/job-detail/9605b724-b3bc-4eae-993a-e972463ab331", "Alpharetta, Georgia", "Contract", "USD 87.00 - 89.43 per hour", "SolutionIT, Inc.")
)

$startRow = 77
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
